$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G (so old G,H shift to H,I)
$ws.Range("G1").EntireColumn.Insert()

# New header for inserted column G1
$ws.Range("G1").Value = "d=6"

# New values for column G (d=6 data), rows 2-10
$ws.Range("G2").Value = 984459420.3467522
$ws.Range("G3").Value = 279671436.6734123
$ws.Range("G4").Value = 1017154738.478231
$ws.Range("G5").Value = 919666473.7647499
$ws.Range("G6").Value = 261122965.848243
$ws.Range("G7").Value = 29188227.75817989
$ws.Range("G8").Value = 186453275.2533088
$ws.Range("G9").Value = 528262424.9802983
$ws.Range("G10").Value = 29.52513657142902
